# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (populated with fresh fund-holding
# data) right after the "总计" summary sheet, pushing the existing
# "2022-Q1" / "2021-Q4" sheets one slot to the right, and updates the
# "总计" summary sheet with a new top row for 2022-Q3 (shifting its old
# 2022-Q1 / 2021-Q4 rows down by one).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new 2022-Q3 row at the top of the
#    data, pushing the existing two rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give the new row 4 the same formatting as row 3 (border/alignment/etc.)
$summary.Range("A3:D3").Copy()
$summary.Range("A4:D4").PasteSpecial($xlPasteFormats)

# Row 4 <- old row 3 data (2021-Q4), unchanged values
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "'2021-Q4"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.18

# Row 3 <- old row 2 data (2022-Q1), unchanged values
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "'2022-Q1"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.32

# Row 2 <- brand-new 2022-Q3 data
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "'2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.5

# Setting text via a leading apostrophe stamps a stray quote-prefix style
# on column B; restrip it back to the plain (unstyled) look used by the
# rest of that column using a clean, already-unstyled cell as the format
# donor.
$summary.Range("C2").Copy()
$summary.Range("B2").PasteSpecial($xlPasteFormats)
$summary.Range("B3").PasteSpecial($xlPasteFormats)
$summary.Range("B4").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2. New "2022-Q3" sheet: duplicate the "2022-Q1" sheet (so it inherits
#    identical formatting/styles) and drop it right after "总计"; then
#    overwrite its cells with the 2022-Q3 fund data.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(2)   # "2022-Q1"
$q1.Copy($null, $summary)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# --- Row 2 (edit the single row that was copied over) ---
$q3.Range("C2").Value = "'前海开源新兴产业混合A"
$q3.Range("B2").Copy()
$q3.Range("C2").PasteSpecial($xlPasteFormats)

$q3.Range("D2").Value = "'7.73"
$q3.Range("B2").Copy()
$q3.Range("D2").PasteSpecial($xlPasteFormats)

$q3.Range("E2").Value = "'93.97"
$q3.Range("B2").Copy()
$q3.Range("E2").PasteSpecial($xlPasteFormats)

$q3.Range("F2").Value = "'5.26"
$q3.Range("B2").Copy()
$q3.Range("F2").PasteSpecial($xlPasteFormats)

$q3.Range("G2").Value = "'0.4066"
$q3.Range("B2").Copy()
$q3.Range("G2").PasteSpecial($xlPasteFormats)

$q3.Range("H2").Value = 7

# --- Rows 3-5: brand new funds, built off row 2's formatting ---
$q3.Range("A2:H2").Copy()
$q3.Range("A3:H5").PasteSpecial($xlPasteFormats)

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'003300"
$q3.Range("C3").Value = "'华夏圆和灵活配置混合A"
$q3.Range("D3").Value = "'0.77"
$q3.Range("E3").Value = "'75.31"
$q3.Range("F3").Value = "'5.54"
$q3.Range("G3").Value = "'0.0427"
$q3.Range("H3").Value = 10

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'014729"
$q3.Range("C4").Value = "'前海开源新兴产业混合C"
$q3.Range("D4").Value = "'0.60"
$q3.Range("E4").Value = "'93.97"
$q3.Range("F4").Value = "'5.26"
$q3.Range("G4").Value = "'0.0316"
$q3.Range("H4").Value = 7

# Row 5
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'015068"
$q3.Range("C5").Value = "'华夏圆和灵活配置混合C"
$q3.Range("D5").Value = "'0.33"
$q3.Range("E5").Value = "'75.31"
$q3.Range("F5").Value = "'5.54"
$q3.Range("G5").Value = "'0.0183"
$q3.Range("H5").Value = 10

# Re-strip the quote-prefix styling the apostrophe assignments above
# stamped onto B3:G5 so they match the plain (unstyled) look of B2:G2.
$q3.Range("B2:G2").Copy()
$q3.Range("B3:G3").PasteSpecial($xlPasteFormats)
$q3.Range("B4:G4").PasteSpecial($xlPasteFormats)
$q3.Range("B5:G5").PasteSpecial($xlPasteFormats)

# Keep "总计" as the active sheet/tab, matching the un-touched bookViews.
$summary.Activate()
